# Collapse the split "<id>...</id>" runs into a single run per occurrence,
# e.g. "<id>" + "p090r_a1" + "</id>"  ->  "<id>p090r_1</id>"
# Repeated for the three occurrences in the document (a1/a2/a3 -> 1/2/3).

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p090r_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p090r_1</id>", 2)
$d.Content.Find.Execute("<id>p090r_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p090r_2</id>", 2)
$d.Content.Find.Execute("<id>p090r_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p090r_3</id>", 2)
